$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1758241758241758
$ws.Range("C2").Value = 0.608058608058608
$ws.Range("J2").Value = 0.01465201465201465
$ws.Range("P2").Value = 0.1135531135531136
$ws.Range("S2").Value = 0.08791208791208792
# Row 3
$ws.Range("B3").Value = 0.01169590643274854
$ws.Range("C3").Value = 0.03508771929824561
$ws.Range("J3").Value = 0.02923976608187134
$ws.Range("P3").Value = 0.8011695906432749
$ws.Range("S3").Value = 0.1228070175438596
# Row 4
$ws.Range("J4").Value = 0.04081632653061224
$ws.Range("P4").Value = 0.7346938775510204
$ws.Range("S4").Value = 0.2244897959183673
# Row 6
$ws.Range("B6").Value = 0.05490196078431372
$ws.Range("D6").Value = 0.01176470588235294
$ws.Range("F6").Value = 0.06274509803921569
$ws.Range("J6").Value = 0.2196078431372549
$ws.Range("O6").Value = 0.0392156862745098
$ws.Range("Q6").Value = 0.1843137254901961
$ws.Range("R6").Value = 0.06274509803921569
$ws.Range("S6").Value = 0.3647058823529412
# Row 7
$ws.Range("B7").Value = 0.1058823529411765
$ws.Range("D7").Value = 0.01764705882352941
$ws.Range("F7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.1117647058823529
$ws.Range("O7").Value = 0.01176470588235294
$ws.Range("Q7").Value = 0.2176470588235294
$ws.Range("R7").Value = 0.05882352941176471
$ws.Range("S7").Value = 0.4176470588235294
# Row 8
$ws.Range("B8").Value = 0.106508875739645
$ws.Range("D8").Value = 0.01380670611439842
$ws.Range("F8").Value = 0.07297830374753451
$ws.Range("J8").Value = 0.1104536489151874
$ws.Range("O8").Value = 0.01972386587771203
$ws.Range("Q8").Value = 0.1479289940828402
$ws.Range("R8").Value = 0.1124260355029586
$ws.Range("S8").Value = 0.4161735700197239
# Row 9
$ws.Range("B9").Value = 0.09745762711864407
$ws.Range("D9").Value = 0.01694915254237288
$ws.Range("F9").Value = 0.1016949152542373
$ws.Range("J9").Value = 0.09745762711864407
$ws.Range("O9").Value = 0.02542372881355932
$ws.Range("Q9").Value = 0.173728813559322
$ws.Range("R9").Value = 0.1059322033898305
$ws.Range("S9").Value = 0.3813559322033898
# Row 10
$ws.Range("B10").Value = 0.08521870286576169
$ws.Range("D10").Value = 0.0248868778280543
$ws.Range("F10").Value = 0.07239819004524888
$ws.Range("J10").Value = 0.1206636500754148
$ws.Range("O10").Value = 0.01357466063348416
$ws.Range("Q10").Value = 0.220211161387632
$ws.Range("R10").Value = 0.1206636500754148
$ws.Range("S10").Value = 0.3423831070889894
# Row 11
$ws.Range("G11").Value = 0.1870229007633588
$ws.Range("J11").Value = 0.07633587786259542
$ws.Range("K11").Value = 0.2480916030534351
$ws.Range("L11").Value = 0.4580152671755725
$ws.Range("S11").Value = 0.03053435114503817
# Row 12
$ws.Range("G12").Value = 0.7868852459016393
$ws.Range("J12").Value = 0.139344262295082
$ws.Range("K12").Value = 0.00819672131147541
$ws.Range("L12").Value = 0.02459016393442623
$ws.Range("S12").Value = 0.04098360655737705
# Row 13
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2142857142857143
$ws.Range("S13").Value = 0.07142857142857142
# Row 15
$ws.Range("F15").Value = 0.0131578947368421
$ws.Range("H15").Value = 0.1710526315789474
$ws.Range("I15").Value = 0.08771929824561403
$ws.Range("J15").Value = 0.3640350877192983
$ws.Range("K15").Value = 0.05263157894736842
$ws.Range("M15").Value = 0.008771929824561403
$ws.Range("O15").Value = 0.07017543859649122
$ws.Range("S15").Value = 0.2324561403508772
# Row 16
$ws.Range("F16").Value = 0.01
$ws.Range("H16").Value = 0.16
$ws.Range("I16").Value = 0.075
$ws.Range("J16").Value = 0.435
$ws.Range("K16").Value = 0.075
$ws.Range("M16").Value = 0.055
$ws.Range("N16").Value = 0.005
$ws.Range("O16").Value = 0.06
$ws.Range("S16").Value = 0.125
# Row 17
$ws.Range("F17").Value = 0.02040816326530612
$ws.Range("H17").Value = 0.1714285714285714
$ws.Range("I17").Value = 0.1061224489795918
$ws.Range("J17").Value = 0.4306122448979592
$ws.Range("K17").Value = 0.06326530612244897
$ws.Range("M17").Value = 0.01224489795918367
$ws.Range("O17").Value = 0.05510204081632653
$ws.Range("S17").Value = 0.1408163265306122
# Row 18
$ws.Range("F18").Value = 0.01123595505617977
$ws.Range("H18").Value = 0.2134831460674157
$ws.Range("I18").Value = 0.08239700374531835
$ws.Range("J18").Value = 0.4382022471910113
$ws.Range("K18").Value = 0.07116104868913857
$ws.Range("M18").Value = 0.0149812734082397
$ws.Range("O18").Value = 0.0599250936329588
$ws.Range("S18").Value = 0.1086142322097378
# Row 19
$ws.Range("F19").Value = 0.01595744680851064
$ws.Range("H19").Value = 0.2272036474164134
$ws.Range("I19").Value = 0.0972644376899696
$ws.Range("J19").Value = 0.3624620060790273
$ws.Range("K19").Value = 0.09042553191489362
$ws.Range("M19").Value = 0.01595744680851064
$ws.Range("N19").Value = 0.0007598784194528875
$ws.Range("O19").Value = 0.06610942249240122
$ws.Range("S19").Value = 0.1238601823708207
